# Adds three new reference rows (14-16) to Sheet1, each with a title in
# column A and a hyperlinked URL in column B, matching the style already
# used by the existing rows (e.g. row 13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: "A Brief History of Scaling LinkedIn" ---
$ws.Range("A14").Value() = "A Brief History of Scaling LinkedIn"
$ws.Range("B14").Value() = "https://engineering.linkedin.com/architecture/brief-history-scaling-linkedin"
$ws.Hyperlinks.Add($ws.Range("B14"), "https://engineering.linkedin.com/architecture/brief-history-scaling-linkedin")
$ws.Range("B14").Style() = "Hyperlink"
$ws.Rows.Item(14).RowHeight() = 17

# --- Row 15: "Optimizing Linux Memory Management for Low-latency / High-throughput Databases" ---
$ws.Range("B15").Value() = "https://engineering.linkedin.com/performance/optimizing-linux-memory-management-low-latency-high-throughput-databases"
$ws.Range("A15").Value() = "Optimizing Linux Memory Management for Low-latency / High-throughput Databases"
$ws.Hyperlinks.Add($ws.Range("B15"), "https://engineering.linkedin.com/performance/optimizing-linux-memory-management-low-latency-high-throughput-databases")
$ws.Range("B15").Style() = "Hyperlink"
$ws.Rows.Item(15).RowHeight() = 34

# --- Row 16: "Using set cover algorithm to optimize query latency for a large scale distributed graph" ---
$ws.Range("B16").Value() = "https://engineering.linkedin.com/real-time-distributed-graph/using-set-cover-algorithm-optimize-query-latency-large-scale-distributed"
$ws.Range("A16").Value() = "Using set cover algorithm to optimize query latency for a large scale distributed graph"
$ws.Hyperlinks.Add($ws.Range("B16"), "https://engineering.linkedin.com/real-time-distributed-graph/using-set-cover-algorithm-optimize-query-latency-large-scale-distributed")
$ws.Range("B16").Style() = "Hyperlink"
$ws.Rows.Item(16).RowHeight() = 34

# Match the final selection left behind by the author
$null = $ws.Range("A16").Select()
